$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.88572758436203
$ws.Range("B1").Value = 3.190121650695801
$ws.Range("C1").Value = 2.881486654281616
$ws.Range("D1").Value = 1.636927962303162
$ws.Range("E1").Value = 1.257978200912476
